# Adding "Datos normalizados" / "Desempeño" analysis tables (geometric-mean
# normalization) below the two existing tables.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Column A a bit wider to fit the new labels.
# (engine quantizes ColumnWidth to ~1/6-character steps; 16.5 is the
# closest input to the target rendered width)
$ws.Columns.Item(1).ColumnWidth = 16.5

# ---------------------------------------------------------------------
# Section 1: "Datos normalizados" (rows 14-20), referenced off table 1
# (rows 2-6, Pc/computadora).

$ws.Range("A14").Value = "Datos normalizados"
$ws.Range("B14").Value = "Referencia: A"
$ws.Range("A15").Value = "Tiempo de respuesta"

# Header row 16 - copy look & feel from the row-2 header.
$ws.Range("A2:I2").Copy()
$ws.Range("A16:I16").PasteSpecial(-4122)
$ws.Range("A16").Value = $ws.Range("A2").Value2
$ws.Range("B16").Value = $ws.Range("B2").Value2
$ws.Range("C16").Value = $ws.Range("C2").Value2
$ws.Range("D16").Value = $ws.Range("D2").Value2
$ws.Range("E16").Value = $ws.Range("E2").Value2
$ws.Range("F16").Value = $ws.Range("F2").Value2
$ws.Range("G16").Value = $ws.Range("G2").Value2
$ws.Range("H16").Value = $ws.Range("H2").Value2
$ws.Range("I16").Value = $ws.Range("I2").Value2
$ws.Range("J16").Value = "Media geométrica"

# Rows 17-20 - copy formatting from rows 3-6, then drop in the
# normalized formulas (value / reference-row-A value) + GEOMEAN.
$ws.Range("A3:I6").Copy()
$ws.Range("A17:I20").PasteSpecial(-4122)

$ws.Range("A17").Value = $ws.Range("A3").Value2
$ws.Range("A18").Value = $ws.Range("A4").Value2
$ws.Range("A19").Value = $ws.Range("A5").Value2
$ws.Range("A20").Value = $ws.Range("A6").Value2

$cols1 = @("B","C","D","E","F","G","H","I")
foreach ($col in $cols1) {
    $ws.Range($col + "17").Formula = "=" + $col + "3/" + $col + "3"
    $ws.Range($col + "18").Formula = "=" + $col + "4/" + $col + "3"
    $ws.Range($col + "19").Formula = "=" + $col + "5/" + $col + "3"
    $ws.Range($col + "20").Formula = "=" + $col + "6/" + $col + "3"
}

$ws.Range("J17").Formula = "=GEOMEAN(B17:I17)"
$ws.Range("J18").Formula = "=GEOMEAN(B18:I18)"
$ws.Range("J19").Formula = "=GEOMEAN(B19:I19)"
$ws.Range("J20").Formula = "=GEOMEAN(B20:I20)"

# ---------------------------------------------------------------------
# Section 2: "Desempeño" (rows 22-27), referenced off table 2
# (rows 8-12, Pc/Tareas).

$ws.Range("A22").Value = "Desempeño"

$ws.Range("A8:F8").Copy()
$ws.Range("A23:F23").PasteSpecial(-4122)
$ws.Range("A23").Value = $ws.Range("A8").Value2
$ws.Range("B23").Value = $ws.Range("B8").Value2
$ws.Range("C23").Value = $ws.Range("C8").Value2
$ws.Range("D23").Value = $ws.Range("D8").Value2
$ws.Range("E23").Value = $ws.Range("E8").Value2
$ws.Range("F23").Value = $ws.Range("F8").Value2
$ws.Range("G23").Value = "Media geométrica"

$ws.Range("A9:F12").Copy()
$ws.Range("A24:F27").PasteSpecial(-4122)

$ws.Range("A24").Value = $ws.Range("A9").Value2
$ws.Range("A25").Value = $ws.Range("A10").Value2
$ws.Range("A26").Value = $ws.Range("A11").Value2
$ws.Range("A27").Value = $ws.Range("A12").Value2

$cols2 = @("B","C","D","E","F")
foreach ($col in $cols2) {
    $ws.Range($col + "24").Formula = "=" + $col + "9/" + $col + "9"
    $ws.Range($col + "25").Formula = "=" + $col + "10/" + $col + "9"
    $ws.Range($col + "26").Formula = "=" + $col + "11/" + $col + "9"
    $ws.Range($col + "27").Formula = "=" + $col + "12/" + $col + "9"
}

$ws.Range("G24").Formula = "=GEOMEAN(B24:F24)"
$ws.Range("G25").Formula = "=GEOMEAN(B25:F25)"
$ws.Range("G26").Formula = "=GEOMEAN(B26:F26)"
$ws.Range("G27").Formula = "=GEOMEAN(B27:F27)"

# ---------------------------------------------------------------------
# Match author's last selection before saving.
$ws.Range("H23").Select()
